# Full logic part of bimestral classes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (7:50)
$ws.Range("C3").Value = "['ELT-2A-Acionamentos', 0, 0, 0]"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "MEC-1A-Circuitos Elétricos"

# Row 4 (8:40)
$ws.Range("C4").Value = "['ELT-2A-Acionamentos', 0, 0, 0]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "MCT-1A-Circuitos Elétricos"

# Row 6 (9:50)
$ws.Range("C6").Value = "MEC-1A-Circuitos Elétricos"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "ELT-1A-Circuitos Elétricos"
$ws.Range("F6").Value = "MCT-1A-Circuitos Elétricos"

# Row 7 (10:40)
$ws.Range("C7").Value = "[0, 'MCT-2A-Acionamentos', 0, 0]"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "ELT-1A-Circuitos Elétricos"
$ws.Range("F7").Value = "-"
